$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 210, pushing the existing
# rows 210-213 down to 212-215 (their content stays identical).
$ws.Rows("210:211").Insert()

# --- New row 210 ---
$ws.Cells.Item(210, 1).Value = 10
$ws.Cells.Item(210, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(210, 3).Value = "La Araucanía"
$ws.Cells.Item(210, 4).Value = 44448
$ws.Cells.Item(210, 5).Value = 9
$ws.Cells.Item(210, 6).Value = 100112008
$ws.Cells.Item(210, 7).Value = "Coliflor"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 5550
$ws.Cells.Item(210, 11).Value = 800
$ws.Cells.Item(210, 12).Value = 900
$ws.Cells.Item(210, 13).Value = 850
$ws.Cells.Item(210, 14).Value = "$/unidad"
$ws.Cells.Item(210, 15).Value = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value = 850
$ws.Cells.Item(210, 17).Value = 1
$ws.Cells.Item(210, 18).Value = "Hortaliza"

# --- New row 211 ---
$ws.Cells.Item(211, 1).Value = 10
$ws.Cells.Item(211, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(211, 3).Value = "La Araucanía"
$ws.Cells.Item(211, 4).Value = 44448
$ws.Cells.Item(211, 5).Value = 9
$ws.Cells.Item(211, 6).Value = 100112008
$ws.Cells.Item(211, 7).Value = "Coliflor"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 1250
$ws.Cells.Item(211, 11).Value = 800
$ws.Cells.Item(211, 12).Value = 800
$ws.Cells.Item(211, 13).Value = 800
$ws.Cells.Item(211, 14).Value = "$/unidad"
$ws.Cells.Item(211, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(211, 16).Value = 800
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"
